$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Sandsynlighed (E), Effekt (F) and Ranking (G = E*F) for each risk row.
# Writing numeric values into cells that previously held the placeholder "?"
# string (or were blank) converts them to numbers, which also removes the
# now-unused "?" shared string from the workbook once no cell references it.

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 21

$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 40

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 64

$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 25

$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 80

$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 36

$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 27

$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 20

$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 40

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 24

$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 30

# Update the active selection to match the saved view state.
$ws.Range("G15").Select()
